$d = $word.ActiveDocument

# Locate the paragraph that ends with "...FCAAPI." so we can insert the
# new paragraphs right after it (and before the following empty paragraph).
$findRng = $d.Content
$found = $findRng.Find.Execute( `
    "Sets aggregation: FCA Contexts scaled objects / attributes from Sets aggregation. FCAAPI.", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the anchor paragraph text."
}

# Build a fresh, zero-length Range right after the matched text (i.e. right
# before that paragraph's own paragraph mark). Using a freshly constructed
# Range (rather than Collapse()-ing the Find range) is important: it is the
# reliable way to get InsertXML to insert *after* this point instead of
# clobbering the paragraph's own content.
$insertPos = $findRng.End
$rng = $d.Range($insertPos, $insertPos)

# Four new paragraphs, inserted as raw WordprocessingML so the exact run /
# paragraph-mark structure (bold "Deployment:" run followed by a trailing
# empty run, plain rtl runs elsewhere) matches precisely:
#   1) empty paragraph
#   2) "Deployment:" (bold) + trailing empty run
#   3) empty paragraph
#   4) paragraph with the Apache MetaModel / JBoss Teiid ... text
$bodyXml = '<w:p><w:pPr><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p>' + `
    '<w:p><w:pPr><w:rPr/></w:pPr><w:r><w:rPr><w:b w:val="1"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Deployment:</w:t></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p>' + `
    '<w:p><w:pPr><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p>' + `
    '<w:p><w:pPr><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Apache MetaModel. JBoss Teiid. Connectors (I/O). APIs: Model Services (reify data, schema, behavior alignment in Connectors data structures). OpenRefine Knowledge (data, schema, behavior) alignment extensions (Model Services APIs).</w:t></w:r></w:p>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + $bodyXml + '</w:body>' + `
    '</w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'

$rng.InsertXML($xml) | Out-Null

Write-Host "Inserted 4 paragraphs after the FCAAPI. paragraph. Paragraphs.Count =" $d.Paragraphs.Count
